$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meal Calendar")

# --- Week 1 totals row (row 23) ---
# Move the "Total" label from column A to column B (left table)
$ws.Range("A23").ClearContents() | Out-Null
$ws.Range("B23").Value = "Total"

# Move the "Total" label from column G to column H (right table)
$ws.Range("G23").ClearContents() | Out-Null
$ws.Range("H23").Value = "Total"

# Replace the per-row VLOOKUP lookups with real SUM totals over the week's rows
$ws.Range("E23").Formula = "=SUM(E2:E22)"
$ws.Range("F23").Formula = "=SUM(F2:F22)"
$ws.Range("K23").Formula = "=SUM(K2:K22)"
$ws.Range("L23").Formula = "=SUM(L2:L22)"

# --- Week 3 totals row (row 45) ---
$ws.Range("A45").ClearContents() | Out-Null
$ws.Range("B45").Value = "Total"

$ws.Range("G45").ClearContents() | Out-Null
$ws.Range("H45").Value = "Total"

$ws.Range("E45").Formula = "=SUM(E24:E44)"
$ws.Range("F45").Formula = "=SUM(F24:F44)"
$ws.Range("K45").Formula = "=SUM(K24:K44)"
$ws.Range("L45").Formula = "=SUM(L24:L44)"

# Leave the selection where the user ended up after adding the totals
$ws.Activate() | Out-Null
$ws.Range("L46").Select() | Out-Null
